$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02, 1.033724118193505, 1.053372488444328, 1.04432869495017, 1.058495503476715, 1.040484152365476, 1.038847275504187, 1.056118638277352, 1.047100052824963, 1.061227580791977, 1.040322557792851)
    3 = @(1.02, 1.034787669137142, 1.05401881545104, 1.045222706549989, 1.059353671967494, 1.040612638187751, 1.039553081922015, 1.056578424074054, 1.047805064088466, 1.061899681753507, 1.041029366536749)
    4 = @(1.02, 1.035476523672629, 1.054436804874988, 1.045801945050375, 1.059909300257495, 1.040694390328803, 1.040009889982377, 1.05687502181404, 1.048261371768856, 1.062334267414406, 1.041486823316953)
    5 = @(1.02, 1.035766277249511, 1.054612471512971, 1.046045636175057, 1.060142965369367, 1.040728426300761, 1.040201956396754, 1.056999491316653, 1.048453230982708, 1.062516892145063, 1.041679162487669)
    6 = @(1.02, 1.035814937472915, 1.054641963374757, 1.046086563468182, 1.060182203359453, 1.040734121566393, 1.040234206600233, 1.057020377344837, 1.048485446574983, 1.062547551203204, 1.041711458490137)
    7 = @(1.02, 1.035480394749196, 1.054439152363832, 1.045805200560079, 1.059912422193635, 1.040694846426921, 1.040012456288476, 1.056876685848426, 1.048263935294205, 1.062336707951531, 1.041489393267502)
    8 = @(1.02, 1.034083412140533, 1.053590963749656, 1.044630673759313, 1.058785455125227, 1.040527861513099, 1.039085784157855, 1.056274213610524, 1.047338289689924, 1.061454783635289, 1.040561405156164)
    9 = @(1.02, 1.031626863934725, 1.052094679841397, 1.042566822994187, 1.056802230925094, 1.040223019508431, 1.037453694077083, 1.055205633357185, 1.045708126714201, 1.05989840969255, 1.038926997320274)
    10 = @(1.02, 1.02999261884139, 1.051096135790948, 1.04119489443852, 1.055481937888665, 1.040012706417362, 1.036366223903598, 1.054488650761594, 1.044622034370938, 1.058859343992745, 1.037837982814409)
    11 = @(1.02, 1.02928579440466, 1.050663533743051, 1.0406017890041, 1.054910694218656, 1.039919965782448, 1.03589548428482, 1.054177114994729, 1.044151916377906, 1.058409079876097, 1.037366574691377)
    12 = @(1.02, 1.029023370819217, 1.050502813802406, 1.040381626747315, 1.054698578154315, 1.039885266926018, 1.035720652501391, 1.054061236078093, 1.043977319496925, 1.058241781515362, 1.037191494626762)
    13 = @(1.02, 1.029079656024614, 1.050537290209347, 1.040428845788888, 1.054744074613556, 1.039892721290679, 1.03575815351442, 1.054086099749612, 1.044014769937526, 1.058277669836005, 1.037229048895532)
    14 = @(1.02, 1.029264099889933, 1.05065024923152, 1.040583587390048, 1.054893159217603, 1.039917102678079, 1.035881032185917, 1.054167539677671, 1.044137483613044, 1.058395251959326, 1.037352102068836)
    15 = @(1.02, 1.0293777580764, 1.050719842785501, 1.040678947843402, 1.054985024382455, 1.039932091632595, 1.035956744704475, 1.05421769625275, 1.044213095009856, 1.058467691565487, 1.037427922107855)
    16 = @(1.02, 1.030039545621161, 1.051124841546653, 1.041234276964208, 1.055519859088458, 1.04001882609626, 1.03639746838046, 1.054509303766394, 1.044653238140517, 1.058889219434451, 1.037869271662012)
    17 = @(1.02, 1.030454885601384, 1.051378827266076, 1.041582874931625, 1.05585546873243, 1.040072784663383, 1.03667396111729, 1.054691933829186, 1.044929373587218, 1.059153542017832, 1.038146157050241)
    18 = @(1.02, 1.030697224961497, 1.05152695093358, 1.04178629769659, 1.056051267672147, 1.04010409625701, 1.036835248332774, 1.054798354802489, 1.045090454798783, 1.059307683852846, 1.038307673312072)
    19 = @(1.02, 1.030779869745213, 1.051577453510266, 1.041855675087336, 1.056118037389395, 1.040114745285583, 1.036890245409944, 1.054834623917623, 1.045145382033115, 1.059360236593861, 1.038362748491399)
    20 = @(1.02, 1.030410315423767, 1.051351579255868, 1.041545464200686, 1.055819456483159, 1.040067012122832, 1.036644294644572, 1.054672350105556, 1.044899745198618, 1.059125186117058, 1.038116448447729)
    21 = @(1.02, 1.029209782386674, 1.050616986492672, 1.040538015875893, 1.054849255594359, 1.039909929887231, 1.035844846871942, 1.05414356207387, 1.044101346773814, 1.058360628331877, 1.037315865367564)
    22 = @(1.02, 1.028455668321655, 1.050154932930514, 1.039905423916948, 1.054239653237856, 1.039809714974176, 1.035342329121391, 1.05381016311918, 1.043599512069873, 1.057879630511471, 1.036812633984168)
    23 = @(1.02, 1.028855371238113, 1.050399893338421, 1.040240693806598, 1.054562776470422, 1.039862978142166, 1.035608711097512, 1.053986991854823, 1.043865529662088, 1.05813464366679, 1.037079394253249)
    24 = @(1.02, 1.030430454525636, 1.051363891526659, 1.041562368217106, 1.05583572873414, 1.04006962098457, 1.036657699604336, 1.05468119946891, 1.04491313294244, 1.05913799903027, 1.038129872444074)
    25 = @(1.02, 1.032261332954627, 1.052481692399356, 1.043099682339992, 1.057314620876504, 1.040303079991812, 1.037875527118313, 1.055482702490929, 1.046129446135158, 1.060301036219514, 1.039349429412807)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $cellRef = "$($cols[$i])$r"
        $ws.Range($cellRef).Value = $rowVals[$i]
    }
}